# ---------------------------------------------------------------------------
# Update the validation plots/results for sample 20220721-24001.
# Re-running the peak-calling pipeline lowered several measured peak heights
# (w_height / m_height) below their min_height threshold for a handful of
# markers, which flips their detection status, fills in genotype/phenotype
# calls on marker_table, and produces a final diplotype call on
# genotype_result.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$peak    = $wb.Worksheets.Item("peak_table")
$allele  = $wb.Worksheets.Item("allele_table")
$marker  = $wb.Worksheets.Item("marker_table")
$geno    = $wb.Worksheets.Item("genotype_result")

# ---------------------------------------------------------------------------
# 1) peak_table: lowered w_height (col N) / m_height (col O) values
# ---------------------------------------------------------------------------
$peak.Cells.Item(2, 14).Value = 900    # S1 CYP2D6_14  w_height
$peak.Cells.Item(3, 14).Value = 700    # S1 CYP2D6_10B w_height
$peak.Cells.Item(3, 15).Value = 700    # S1 CYP2D6_10B m_height
$peak.Cells.Item(4, 14).Value = 500    # S1 CYP2D6_49  w_height
$peak.Cells.Item(6, 14).Value = 500    # S1 CYP2D6_41  w_height
$peak.Cells.Item(12, 14).Value = 800   # S2 CYP2D6_4   w_height

# ---------------------------------------------------------------------------
# 2) allele_table: corresponding min_height / is_detected / peak call columns
#    K = min_height, M = is_detected, N = peak, O = size, P = height,
#    Q = status, R = message
# ---------------------------------------------------------------------------

function Set-AlleleDetected($row, $minHeight, $peakNo, $size, $height) {
    $allele.Cells.Item($row, 11).Value = $minHeight
    $allele.Cells.Item($row, 13).Value = $true
    $allele.Cells.Item($row, 14).Value = $peakNo
    $allele.Cells.Item($row, 15).Value = $size
    $allele.Cells.Item($row, 16).Value = $height
    $allele.Cells.Item($row, 17).Value = "ok"
    $allele.Cells.Item($row, 18).Value = ""
}

Set-AlleleDetected 2  900 38 29.22 943    # S1 CYP2D6_14  wildtype G
Set-AlleleDetected 4  700 42 32.59 711    # S1 CYP2D6_10B wildtype C
Set-AlleleDetected 5  700 38 35.16 740    # S1 CYP2D6_10B mutant   T
Set-AlleleDetected 6  500 17 38.87 508    # S1 CYP2D6_49  wildtype T
Set-AlleleDetected 10 500 18 46.91 706    # S1 CYP2D6_41  wildtype G
Set-AlleleDetected 22 800 42 30.64 880    # S2 CYP2D6_4   wildtype G

# ---------------------------------------------------------------------------
# 3) marker_table: genotype (G) / phenotype (H) calls now resolvable
# ---------------------------------------------------------------------------

$marker.Cells.Item(2, 7).Value = "GG"
$marker.Cells.Item(2, 8).Value = "wildtype"

$marker.Cells.Item(3, 7).Value = "CT"
$marker.Cells.Item(3, 8).Value = "heterozygous"

$marker.Cells.Item(4, 7).Value = "TT"
$marker.Cells.Item(4, 8).Value = "wildtype"

$marker.Cells.Item(6, 7).Value = "GG"
$marker.Cells.Item(6, 8).Value = "wildtype"

$marker.Cells.Item(12, 7).Value = "GG"
$marker.Cells.Item(12, 8).Value = "wildtype"

# ---------------------------------------------------------------------------
# 4) genotype_result: final diplotype call for the sample
# ---------------------------------------------------------------------------

$geno.Cells.Item(2, 2).Value = "*1/*10B"
